$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# Range.InsertXML replaces the whole paragraph(s) touched by the target
# range, so each replacement fragment re-specifies the paragraph's own
# attributes (paraId/textId/rsid*) together with its original <w:pPr> so
# those stay exactly as they were before the edit — only the run/text
# content inside the paragraph actually changes.

# --- Row 7 ("-Frais de facturation"): quantity cell -> "caca" ----------
$qtyCell = $t.Cell(7, 2)
$qtyPara = $qtyCell.Range.Paragraphs.Item(1).Range
$qtyXml = '<w:p ' + $wNs + ' w14:paraId="6DC44C79" w14:textId="03E48F47" w:rsidR="00900638" w:rsidRPr="009F2189" w:rsidRDefault="00E11583" w:rsidP="00900638"><w:pPr><w:ind w:firstLine="243"/></w:pPr><w:r><w:t>caca</w:t></w:r></w:p>'
$null = $qtyPara.InsertXML($qtyXml)

# --- Row 7: unit-price cell "       5,00" -> "   " + "1040" + ",00" ----
$priceCell = $t.Cell(7, 3)
$pricePara = $priceCell.Range.Paragraphs.Item(1).Range
$priceXml = '<w:p ' + $wNs + ' w14:paraId="1EB6AD4E" w14:textId="5DF3CA76" w:rsidR="00900638" w:rsidRPr="009F2189" w:rsidRDefault="00900638" w:rsidP="00900638"><w:pPr><w:ind w:firstLine="243"/></w:pPr><w:r><w:t xml:space="preserve">   </w:t></w:r><w:r><w:t>1040</w:t></w:r><w:r><w:t>,00</w:t></w:r></w:p>'
$null = $pricePara.InsertXML($priceXml)

# --- Row 23 ("-Fourniture et mise en place de béton..."): "3,00" -> "3" + "." + "00" ----
$qty2Cell = $t.Cell(23, 2)
$qty2Para = $qty2Cell.Range.Paragraphs.Item(1).Range
$qty2Xml = '<w:p ' + $wNs + ' w14:paraId="0EA6F3DE" w14:textId="0AEF39B8" w:rsidR="00D2409B" w:rsidRPr="009F2189" w:rsidRDefault="00E11583" w:rsidP="00D2409B"><w:pPr><w:ind w:firstLine="243"/></w:pPr><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:t>3</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>00</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>m&#179;</w:t></w:r></w:p>'
$null = $qty2Para.InsertXML($qty2Xml)

Write-Host "Edits applied."
